$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters -> indices: D=4, M=13, N=14, O=15, P=16, S=19
# Apply new values per row as derived from the target diff (rows 2-20 data
# shuffle: each row's Fecha/Volumen/Precio*/Precio-$-Kg tuple is replaced
# with another row's tuple).

$ws.Cells.Item(2, 4).Value = 44400

$ws.Cells.Item(3, 4).Value = 44307
$ws.Cells.Item(3, 13).Value = 30

$ws.Cells.Item(4, 4).Value = 44301
$ws.Cells.Item(4, 13).Value = 38
$ws.Cells.Item(4, 14).Value = 22000
$ws.Cells.Item(4, 15).Value = 22000
$ws.Cells.Item(4, 16).Value = 22000
$ws.Cells.Item(4, 19).Value = 1100

$ws.Cells.Item(5, 4).Value = 44305
$ws.Cells.Item(5, 13).Value = 20

$ws.Cells.Item(6, 4).Value = 44377
$ws.Cells.Item(6, 14).Value = 20000
$ws.Cells.Item(6, 15).Value = 20000
$ws.Cells.Item(6, 16).Value = 20000
$ws.Cells.Item(6, 19).Value = 1000

$ws.Cells.Item(7, 4).Value = 44445
$ws.Cells.Item(7, 13).Value = 45

$ws.Cells.Item(8, 4).Value = 44448
$ws.Cells.Item(8, 13).Value = 30
$ws.Cells.Item(8, 14).Value = 22000
$ws.Cells.Item(8, 15).Value = 22000
$ws.Cells.Item(8, 16).Value = 22000
$ws.Cells.Item(8, 19).Value = 1100

$ws.Cells.Item(9, 4).Value = 44406
$ws.Cells.Item(9, 13).Value = 20
$ws.Cells.Item(9, 14).Value = 20000
$ws.Cells.Item(9, 15).Value = 20000
$ws.Cells.Item(9, 16).Value = 20000
$ws.Cells.Item(9, 19).Value = 1000

$ws.Cells.Item(10, 4).Value = 44300
$ws.Cells.Item(10, 13).Value = 45
$ws.Cells.Item(10, 14).Value = 22000
$ws.Cells.Item(10, 15).Value = 22000
$ws.Cells.Item(10, 16).Value = 22000
$ws.Cells.Item(10, 19).Value = 1100

$ws.Cells.Item(11, 4).Value = 44382
$ws.Cells.Item(11, 13).Value = 24
$ws.Cells.Item(11, 14).Value = 20000
$ws.Cells.Item(11, 15).Value = 20000
$ws.Cells.Item(11, 16).Value = 20000
$ws.Cells.Item(11, 19).Value = 1000

$ws.Cells.Item(12, 4).Value = 44294
$ws.Cells.Item(12, 13).Value = 25

$ws.Cells.Item(13, 4).Value = 44376
$ws.Cells.Item(13, 13).Value = 38

$ws.Cells.Item(14, 4).Value = 44292
$ws.Cells.Item(14, 13).Value = 30
$ws.Cells.Item(14, 14).Value = 25000
$ws.Cells.Item(14, 15).Value = 25000
$ws.Cells.Item(14, 16).Value = 25000
$ws.Cells.Item(14, 19).Value = 1250

$ws.Cells.Item(16, 4).Value = 44291
$ws.Cells.Item(16, 13).Value = 70
$ws.Cells.Item(16, 14).Value = 25000
$ws.Cells.Item(16, 15).Value = 25000
$ws.Cells.Item(16, 16).Value = 25000
$ws.Cells.Item(16, 19).Value = 1250

$ws.Cells.Item(17, 4).Value = 44389
$ws.Cells.Item(17, 13).Value = 20

$ws.Cells.Item(18, 4).Value = 44385
$ws.Cells.Item(18, 13).Value = 36

$ws.Cells.Item(19, 4).Value = 44413
$ws.Cells.Item(19, 14).Value = 20000
$ws.Cells.Item(19, 15).Value = 20000
$ws.Cells.Item(19, 16).Value = 20000
$ws.Cells.Item(19, 19).Value = 1000

$ws.Cells.Item(20, 4).Value = 44298
$ws.Cells.Item(20, 13).Value = 65
